$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    "B4" = 0.6718380303341256
    "C4" = 0.6809999999999999
    "D4" = 0.6855582855486929
    "E4" = 0.677
    "F4" = 0.5907062303113287
    "G4" = 0.597
    "H4" = 0.5916427939661594
    "I4" = 0.589
    "J4" = 0.9287481116485458
    "K4" = 0.93
    "L4" = 0.9329904203633393
    "M4" = 0.9279999999999999

    "B5" = 0.6800333987524947
    "C5" = 0.6980000000000001
    "D5" = 0.6827706965798201
    "E5" = 0.6785000000000001
    "F5" = 0.5807759710969893
    "G5" = 0.576
    "H5" = 0.6303047285603531
    "I5" = 0.5990000000000001
    "J5" = 0.8129093776788885
    "K5" = 0.8160000000000001
    "L5" = 0.8441411950493343
    "M5" = 0.8179999999999999

    "B6" = 0.7008501192909451
    "C6" = 0.8160000000000001
    "D6" = 0.6208081766012383
    "E6" = 0.6530000000000001
    "F6" = 0.5963788963570399
    "G6" = 0.603
    "H6" = 0.5969462172508369
    "I6" = 0.5940000000000001
    "J6" = 0.7512941382518582
    "K6" = 0.729
    "L6" = 0.819597667191821
    "M6" = 0.7675

    "B7" = 0.4400990361462195
    "C7" = 0.441
    "D7" = 0.4447796887164944
    "E7" = 0.4415
    "F7" = 0.6666666666666667
    "G7" = 1
    "H7" = 0.5
    "I7" = 0.5
    "J7" = 0.8306033025396022
    "K7" = 0.8230000000000001
    "L7" = 0.868980878121641
    "M7" = 0.8375
}

foreach ($addr in $data.Keys) {
    $ws.Range($addr).Value = $data[$addr]
}
